$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: update count (C31) and recompute ratio (E31)
$ws.Range("C31").Value = 38
$ws.Range("E31").Value = 0.01643598615916955

# Row 37: update count (C37) and total (D37)
$ws.Range("C37").Value = 588
$ws.Range("D37").Value = 588
